# Applies the "DECEMBER 21" rent-ledger update (payments recorded for several
# tenants, a new reconciliation block in columns J:L, a "NOV" -> "DECEMBER"
# label fix, and the view/selection state left by the editor).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "DECEMBER 21": record this month's payments (column G = "PAID") for
# several tenants, and correct one tenant's rent-due split (E column).
# ---------------------------------------------------------------------------
$dec = $wb.Worksheets.Item("DECEMBER 21")

# Row 7 - payment of 6000 received (fully settles the balance).
$dec.Range("G7").Value = 6000

# Row 9 - the 2500 previously logged as "RENT" (E9) actually belongs to the
# B/F column; clear it so F9/H9 fall back to the brought-forward amount only.
$dec.Range("E9").ClearContents()

# Row 11 - payment of 2500 received (fully settles the balance).
$dec.Range("G11").Value = 2500

# Row 16 - rent of 2500 now due for this tenant.
$dec.Range("E16").Value = 2500

# Row 17 - partial payment of 1550 received.
$dec.Range("G17").Value = 1550

# Row 19 - payment of 2500 received (fully settles the balance).
$dec.Range("G19").Value = 2500

# Row 20 - partial payment of 2000 received.
$dec.Range("G20").Value = 2000

# Row 29 - the summary label was "NOV"; fix it to read "DECEMBER".
$dec.Range("A29").Value = "DECEMBER"
$dec.Range("E29").Value = "DECEMBER"

# New reconciliation block (columns J:L) added alongside the existing summary.
$dec.Range("J31").Formula = "=B29-C33"
$dec.Range("J31").NumberFormat = "#,##0"

$dec.Range("J33").Formula = "=J31-J32"
$dec.Range("J33").NumberFormat = "#,##0"

$dec.Range("J34").Value = 29900
$dec.Range("K34").Formula = "=J33+J34"
$dec.Range("K34").NumberFormat = "#,##0"

$dec.Range("J35").Formula = "=J33-J34"
$dec.Range("J35").NumberFormat = "#,##0"

$dec.Range("J36").Value = 2500

$dec.Range("J37").Formula = "=J35-J36"
$dec.Range("J37").NumberFormat = "#,##0"

$dec.Range("L41").Formula = "=15800+3000"

# Row 38 - record the 29900 payment made on 11/12 on both halves of the sheet.
$dec.Range("A38").Value = "PAID ON 11/12"
$dec.Range("C38").Value = 29900
$dec.Range("E38").Value = "PAID ON 11/12"
$dec.Range("G38").Value = 29900

# ---------------------------------------------------------------------------
# Leave the selection/scroll state the way the editor left each tab. Touch
# the sheets in tab order, finishing on "DECEMBER 21" (the last-active tab),
# so tabSelected ends up on the right sheet.
# ---------------------------------------------------------------------------
$oct = $wb.Worksheets.Item("OCTOBER 21")
$oct.Range("E16").Select()

$nov = $wb.Worksheets.Item("NOVEMBER 21")
$nov.Range("A20").Select()

$dec.Range("K42").Select()
